{"js": "// Resume edit: add \"Panic \" (non-italic) before the italic \"Visual Studio\n// Code\" run, and change that run's text to \"Nova, Visual Studio Code\",\n// turning the \"Development Software\" list entry into \"Panic Nova, Visual\n// Studio Code, Visual Studio, Xcode, Eclipse, ...\".\n\nconst body = context.document.body;\n\n// Locate the unique \"Visual Studio Code\" run of text in the document.\nconst results = body.search(\"Visual Studio Code\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Visual Studio Code\" in the document.');\n}\n\nconst target = results.items[0];\n\n// Insert the new, non-italic \"Panic \" text immediately before it; it\n// inherits the surrounding (Tahoma/21 half-points) character formatting\n// but must explicitly be non-italic since the \"Visual Studio Code\" run\n// right after it is italic.\nconst panicRange = target.insertText(\"Panic \", Word.InsertLocation.before);\npanicRange.font.set({ italic: false });\n\n// Replace the original run's text in place so its existing (italic)\n// formatting is preserved, turning \"Visual Studio Code\" into\n// \"Nova, Visual Studio Code\".\ntarget.insertText(\"Nova, Visual Studio Code\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Resume edit: add \"Panic \" (non-italic) before the italic \"Visual Studio\n# Code\" run, and change that run's text to \"Nova, Visual Studio Code\",\n# turning the \"Development Software\" list entry into \"Panic Nova, Visual\n# Studio Code, Visual Studio, Xcode, Eclipse, ...\".\n\n$d = $word.ActiveDocument\n\n# Locate the unique \"Visual Studio Code\" text.\n$vscRange = $d.Content\n$find = $vscRange.Find\n$find.Text = \"Visual Studio Code\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"Visual Studio Code\" in the document.'\n}\n\n$matchStart = $vscRange.Start\n$matchLength = $vscRange.End - $vscRange.Start\n\n# Insert the new, non-italic \"Panic \" text immediately before the match.\n# InsertBefore collapses $vscRange onto the newly inserted text, so grab a\n# fresh Range over those bounds to set its formatting explicitly\n# non-italic (it would otherwise inherit the italic run that follows it).\n$vscRange.InsertBefore(\"Panic \")\n$panicRange = $d.Range($matchStart, $matchStart + 6)\n$panicRange.Font.Italic = 0\n\n# Replace the original run's text in place so its existing (italic)\n# formatting is preserved, turning \"Visual Studio Code\" into\n# \"Nova, Visual Studio Code\".\n$vscRange2 = $d.Range($matchStart + 6, $matchStart + 6 + $matchLength)\n$vscRange2.Text = \"Nova, Visual Studio Code\"\n"}
